$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AA2").Value = 0
$ws.Range("AA5").Value = 2
$ws.Range("AA7").Value = 0
$ws.Range("AA8").Value = 0
$ws.Range("AA10").Value = 2
$ws.Range("AA12").Value = 0
$ws.Range("AA13").Value = 0
$ws.Range("AA14").Value = 0
$ws.Range("AA15").Value = 0
$ws.Range("AA16").Value = 2
$ws.Range("AA17").Value = 2
$ws.Range("AA18").Value = 2
$ws.Range("AA19").Value = 2
$ws.Range("AA20").Value = 0
$ws.Range("AA21").Value = 0
$ws.Range("AA22").Value = 0
$ws.Range("AA23").Value = 0
$ws.Range("AA28").Value = 1
$ws.Range("AA29").Value = 2
$ws.Range("AA32").Value = 2
$ws.Range("AA33").Value = 2
$ws.Range("AA35").Value = 2
$ws.Range("AA36").Value = 2
$ws.Range("AA38").Value = 2
$ws.Range("AA39").Value = 2
$ws.Range("AA40").Value = 0
$ws.Range("AA42").Value = 0
$ws.Range("AA43").Value = 0
$ws.Range("AA44").Value = 0
$ws.Range("AA46").Value = 0
$ws.Range("AA48").Value = 0
$ws.Range("AA49").Value = 0
$ws.Range("AA50").Value = 0
$ws.Range("AA51").Value = 2
$ws.Range("AA53").Value = 0
$ws.Range("AA54").Value = 2
$ws.Range("AA55").Value = 2
$ws.Range("AA56").Value = 2
$ws.Range("AA57").Value = 1
$ws.Range("AA63").Value = 1
$ws.Range("AA64").Value = 1
$ws.Range("AA66").Value = 2
$ws.Range("AA67").Value = 0
$ws.Range("AA69").Value = 0
$ws.Range("AA70").Value = 2
$ws.Range("AA71").Value = 2
$ws.Range("AA72").Value = 1
$ws.Range("AA73").Value = 1
$ws.Range("AA76").Value = 1
$ws.Range("AA77").Value = 2
$ws.Range("AA79").Value = 0
$ws.Range("AA80").Value = 2
$ws.Range("AA82").Value = 0
$ws.Range("AA83").Value = 0
$ws.Range("AA84").Value = 2
$ws.Range("AA85").Value = 2
$ws.Range("AA88").Value = 2
$ws.Range("AA89").Value = 0
$ws.Range("AA90").Value = 0
$ws.Range("AA91").Value = 1
$ws.Range("AA92").Value = 2
$ws.Range("AA93").Value = 2
$ws.Range("AA95").Value = 2
$ws.Range("AA97").Value = 2
$ws.Range("AA98").Value = 0
$ws.Range("AA99").Value = 0
$ws.Range("AA102").Value = 0
$ws.Range("AA104").Value = 2
$ws.Range("AA105").Value = 2
$ws.Range("AA109").Value = 2
$ws.Range("AA110").Value = 0
$ws.Range("AA111").Value = 1
$ws.Range("AA112").Value = 0
$ws.Range("AA113").Value = 0
$ws.Range("AA114").Value = 2
$ws.Range("AA119").Value = 0
$ws.Range("AA123").Value = 2
$ws.Range("AA124").Value = 2
$ws.Range("AA126").Value = 0
$ws.Range("AA127").Value = 2
$ws.Range("AA128").Value = 0
$ws.Range("AA129").Value = 0
$ws.Range("AA130").Value = 0
$ws.Range("AA132").Value = 1
$ws.Range("AA133").Value = 0
$ws.Range("AA134").Value = 2
$ws.Range("AA137").Value = 1
$ws.Range("AA138").Value = 2
$ws.Range("AA139").Value = 0
$ws.Range("AA140").Value = 2
$ws.Range("AA141").Value = 2
$ws.Range("AA142").Value = 2
$ws.Range("AA143").Value = 2
$ws.Range("AA145").Value = 1
$ws.Range("AA147").Value = 0
$ws.Range("AA149").Value = 2
$ws.Range("AA150").Value = 2
$ws.Range("AA151").Value = 0
$ws.Range("AA152").Value = 0
$ws.Range("AA153").Value = 2
$ws.Range("AA155").Value = 2
$ws.Range("AA156").Value = 2
$ws.Range("AA158").Value = 1
$ws.Range("AA159").Value = 1
$ws.Range("AA160").Value = 2
